$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (years header): add 2019 (P4) and 2020 (Q4), matching O4's style ---
$ws.Range("O4").Copy($ws.Range("P4"))
$ws.Range("P4").Value = 2019
$ws.Range("O4").Copy($ws.Range("Q4"))
$ws.Range("Q4").Value = 2020

# --- Row 5 (share %): add 35.67 (P5) and a blank but formatted cell (Q5), matching E5's style ---
$ws.Range("E5").Copy($ws.Range("P5"))
$ws.Range("P5").Value = 35.67
$ws.Range("E5").Copy($ws.Range("Q5"))
$ws.Range("Q5").ClearContents()

# --- Row 6 (hydro output): add 13859.3 (P6) and 13979.1 (Q6), matching D6's style ---
$ws.Range("D6").Copy($ws.Range("P6"))
$ws.Range("P6").Value = 13859.3
$ws.Range("D6").Copy($ws.Range("Q6"))
$ws.Range("Q6").Value = 13979.1

# --- Update the view's active selection to P9, as recorded in the saved file ---
$ws.Range("P9").Select()
